$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 133 (existing rows 133:170 shift down to 138:175),
# matching the weekly "Fruta, Vega Modelo de Temuco - Frutilla" refresh.
$ws.Rows("133:137").Insert()

# Columns that are constant across every data row in this sheet.
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = "Fruta"
$constG = 100101
$constH = "Berries"
$constI = 100112025
$constJ = "Frutilla"
$constK = "Sin especificar"
$constT = 7

# New weekly rows (date 44511) to populate at rows 133-137.
$newRows = @(
    @{ Row = 133; D = 44511; L = "Primera"; M = 280;  N = 9000;  O = 9000;  P = 9000;  Q = "`$/bandeja 7 kilos"; R = "Provincia de Melipilla";    S = 1286 },
    @{ Row = 134; D = 44511; L = "Primera"; M = 450;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/bandeja 7 kilos"; R = "Región del Maule";          S = 1143 },
    @{ Row = 135; D = 44511; L = "Primera"; M = 180;  N = 10000; O = 10000; P = 10000; Q = "`$/caja 7 kilos";    R = "Región de La Araucanía";    S = 1429 },
    @{ Row = 136; D = 44511; L = "Segunda"; M = 150;  N = 7000;  O = 7000;  P = 7000;  Q = "`$/bandeja 7 kilos"; R = "Provincia de Melipilla";    S = 1000 },
    @{ Row = 137; D = 44511; L = "Segunda"; M = 110;  N = 7000;  O = 7000;  P = 7000;  Q = "`$/caja 7 kilos";    R = "Región de La Araucanía";    S = 1000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $constI
    $ws.Cells.Item($row, 10).Value = $constJ
    $ws.Cells.Item($row, 11).Value = $constK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $constT
}
